# 14 apr cdc update
# Updates the daily CDC "cases" counts (column B) for the existing date
# range (rows 4-90), un-highlights rows 80-83 (no longer among the most
# recent 4 days of preliminary data) and appends 4 new daily rows
# (91-94, dates 43931-43934) highlighted as the new preliminary days.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Updated case counts for existing rows ---------------------------------
$updates = @{
    4  = 3
    7  = 0
    12 = 3
    13 = 1
    14 = 5
    15 = 5
    16 = 1
    17 = 2
    21 = 4
    22 = 7
    23 = 7
    24 = 3
    25 = 1
    27 = 4
    28 = 4
    30 = 3
    33 = 9
    35 = 14
    36 = 13
    37 = 10
    39 = 27
    41 = 36
    42 = 42
    43 = 34
    44 = 51
    45 = 99
    46 = 80
    47 = 116
    48 = 117
    49 = 181
    50 = 160
    51 = 374
    52 = 366
    53 = 421
    54 = 466
    55 = 547
    56 = 730
    57 = 872
    58 = 1377
    59 = 2325
    60 = 2514
    61 = 3327
    62 = 5205
    63 = 6566
    64 = 7387
    65 = 9042
    66 = 10675
    67 = 8890
    68 = 8733
    69 = 10407
    70 = 12445
    71 = 12412
    72 = 12817
    73 = 14190
    74 = 10272
    75 = 10545
    76 = 15451
    77 = 14728
    78 = 14088
    79 = 14092
    80 = 14651
    81 = 10849
    82 = 10194
    83 = 15451
    84 = 14774
    85 = 12417
    86 = 9328
    87 = 7911
    88 = 3525
    89 = 1232
    90 = 365
}

foreach ($row in $updates.Keys) {
    $ws.Cells.Item($row, 2).Value = $updates[$row]
}

# --- Rows 80-83 are no longer among the newest 4 days of preliminary data,
#     so they lose the yellow "recent" highlight on column A. Copy the
#     (unhighlighted) format from a row that has already rolled off the
#     highlight window rather than clearing the fill by hand, so the
#     existing unhighlighted style is reused instead of a new one minted. --
$formatSource = $ws.Range("A79")
$unhighlightTarget = $ws.Range("A80:A83")
$formatSource.Copy()
$unhighlightTarget.PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Append the 4 new daily rows, highlighted as the newest data --------
$newRows = @{
    91 = @{ Date = 43931; Cases = 144 }
    92 = @{ Date = 43932; Cases = 27 }
    93 = @{ Date = 43933; Cases = 9 }
    94 = @{ Date = 43934; Cases = 0 }
}

foreach ($row in ($newRows.Keys | Sort-Object)) {
    $info = $newRows[$row]
    $dateCell = $ws.Cells.Item($row, 1)
    $dateCell.Value = $info.Date
    $dateCell.Interior.ColorIndex = 6
    $ws.Cells.Item($row, 2).Value = $info.Cases
}
